$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" column (K), mirroring the existing yearly columns (D:J).
# Copy the formatting from the neighboring 2021 column (J) first, so the new
# cells pick up the exact same styles (number format / font / borders) that
# the rest of the row already uses, then fill in the 2022 values.
$ws.Range("J4:J14").Copy()
$ws.Range("K4:K14").PasteSpecial(-4122)

$ws.Range("K4").Value = 2022
$ws.Range("K5").Value = 1.6
$ws.Range("K6").Value = 0.4
$ws.Range("K7").Value = 0.9
$ws.Range("K8").Value = 0.6
$ws.Range("K9").Value = 2.1
$ws.Range("K10").Value = 0.6
$ws.Range("K11").Value = 0.9
$ws.Range("K12").Value = 2.3
$ws.Range("K13").Value = 4.3
$ws.Range("K14").Value = 0.3

# Clear the marching-ants marquee left over from the Copy above.
$excel.CutCopyMode = 0

# Update the saved selection to match the source workbook (active cell L7).
$ws.Range("L7").Select()
